# Apply the "Name" column (B) shuffle described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column B, as per the target diff.
$changes = @{
    2  = "Franck"
    5  = "None"
    7  = "None"
    8  = "arr"
    9  = "None"
    10 = "KEvin"
    11 = "Dick"
    12 = "Tom"
    14 = "None"
    15 = "None"
    18 = "None"
    19 = "None"
    21 = "Evris"
    22 = "Jean"
}

foreach ($row in $changes.Keys) {
    $ws.Cells.Item($row, 2).Value = $changes[$row]
}
